$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet2: rename to the "helper" sheet ---
$ws2.Name = "辅助信息（请勿操作）"

# --- Sheet1: unfreeze panes / collapse multi-pane selection ---
$ws1.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws1.Range("G3").Select()

# --- Sheet2: populate helper content ---
$ws2.Range("A1").Value = "校验码"
$ws2.Range("A2").Value = "H86D$8#a"
$ws2.Range("B1").Value = "日期格式"
$ws2.Range("A1").Copy()
$ws2.Range("B2").PasteSpecial(-4122)
$ws2.Range("B2").NumberFormat = "yyyy\-mm\-dd;@"
$ws2.Range("B2").Value = [datetime]"2022-01-01"
$ws2.Columns.Item(2).ColumnWidth = 10.33203125

# Make Sheet2 the active/selected tab, matching the authored state
$ws2.Activate()
$ws2.Range("C2").Select()
